$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add dropdown selection value "null" to D4 (matching existing C/D column "null" values)
$ws.Range("D4").Value = "null"

# Update the active cell selection to F9
$ws.Range("F9").Select()
